$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name to reflect the new "through" date
$ws.Name = "Through 2022-08-02"

# Update July row (row 8) - only I8 changes
$ws.Range("I8").Value = 166

# Update the label in A9 (shared string used for row label "August (through 08-0X)")
$ws.Range("A9").Value = "August (through 08-02)"

# Update August row (row 9)
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 13
$ws.Range("H9").Value = 10
$ws.Range("I9").Value = 6

# Update Total row (row 10)
$ws.Range("B10").Value = 163
$ws.Range("C10").Value = 305
$ws.Range("D10").Value = 469
$ws.Range("E10").Value = 429
$ws.Range("F10").Value = 309
$ws.Range("G10").Value = 634
$ws.Range("H10").Value = 920
$ws.Range("I10").Value = 977
